$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: fix country order (D2)
$ws.Range("D2").Value = "['Netherlands', 'United Kingdom', 'United States', 'Australia', 'China']"

# Row 5: typo fix lifestock -> livestock (B5)
$ws.Range("B5").Value = "nutrition and industrial livestock alternatives"

# Row 6: fix country order (D6)
$ws.Range("D6").Value = "['Vietnam', 'Sri Lanka', 'Thailand', 'United States']"

# Row 7: fix country order (D7)
$ws.Range("D7").Value = "['United States', 'Brazil', 'China']"

# Row 8: fix country order (D8)
$ws.Range("D8").Value = "['United States', 'India', 'Philippines', 'China']"

# Row 11: fix country order (D11)
$ws.Range("D11").Value = "['Brazil', 'Canada', 'United States', 'Mexico', 'India']"

# Row 12: typo fix lifestock -> livestock (B12)
$ws.Range("B12").Value = "nutrition and industrial livestock alternatives"

# Row 14: fix country order (D14) and continent order (E14)
$ws.Range("D14").Value = "['United Kingdom', 'Vietnam', 'India']"
$ws.Range("E14").Value = "['Asia', 'Europe']"

# Row 15: fix country order (D15)
$ws.Range("D15").Value = "['Uruguay', 'Brazil', 'Thailand', 'Peru', 'Indonesia', 'Ecuador', 'Chile', 'Argentina', 'Austria', 'Colombia']"

# Row 16: typo fix lifestock -> livestock (B16) and fix country order (D16)
$ws.Range("B16").Value = "nutrition and industrial livestock alternatives"
$ws.Range("D16").Value = "['Brazil', 'Malaysia', 'South Korea', 'United Kingdom', 'Thailand', 'United States', 'Indonesia', 'Belgium', 'Singapore', 'Japan', 'India', 'Israel']"

# Row 17: fix country order (D17)
$ws.Range("D17").Value = "['United Kingdom', 'Mexico', 'United States', 'Japan']"
